$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contact/website cell (B10): stat.kg -> stat.gov.kg
$ws.Range("B10").Value = "МФ КР: www.minfin.kg;`nНСК: www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Leave the cursor/selection on the edited cell, as in the source file
$ws.Range("B10").Select() | Out-Null
